$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2125.6191
$ws.Range("I70").Value = 1357.1428
$ws.Range("J70").Value = 2509.8572
$ws.Range("K70").Value = 4071.4284
$ws.Range("L70").Value = 7529.571599999999
$ws.Range("M70").Value = -3801.4284
$ws.Range("N70").Value = -8069.571599999999
$ws.Range("H73").Value = 2125.6191
$ws.Range("I73").Value = 1357.1428
$ws.Range("J73").Value = 2509.8572
$ws.Range("K73").Value = 4071.4284
$ws.Range("L73").Value = 7529.571599999999
$ws.Range("M73").Value = -3135.4284
$ws.Range("N73").Value = -9401.571599999999
$ws.Range("H137").Value = 15626030
$ws.Range("I137").Value = 922.93475
$ws.Range("J137").Value = 55556856
$ws.Range("K137").Value = 2768.80425
$ws.Range("L137").Value = 166670568
$ws.Range("M137").Value = -218.8042500000001
$ws.Range("N137").Value = -166675668
$ws.Range("H138").Value = 2831.016
$ws.Range("I138").Value = 2480.4814
$ws.Range("K138").Value = 7441.4442
$ws.Range("M138").Value = -2301.4442
$ws.Range("H141").Value = 887.449
$ws.Range("I141").Value = 445.97562
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 1337.92686
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = 3842.07314
$ws.Range("N141").Value = -19810
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8363.5
$ws.Range("J37").Value = 10717.167
$ws.Range("L37").Value = 10717.167
$ws.Range("N37").Value = -11263.167
$ws.Range("H55").Value = 17857
$ws.Range("J55").Value = 17857
$ws.Range("L55").Value = 17857
$ws.Range("N55").Value = -18487
$ws.Range("H61").Value = 1985158
$ws.Range("I61").Value = 2268598
$ws.Range("J61").Value = 1077.4286
$ws.Range("K61").Value = 2268598
$ws.Range("L61").Value = 1077.4286
$ws.Range("M61").Value = -2268386
$ws.Range("N61").Value = -1501.4286
$ws.Range("H80").Value = 19975.5
$ws.Range("J80").Value = 19975.5
$ws.Range("L80").Value = 19975.5
$ws.Range("N80").Value = -21971.5
$ws.Range("H83").Value = 19975.5
$ws.Range("J83").Value = 19975.5
$ws.Range("L83").Value = 59926.5
$ws.Range("N83").Value = -69910.5
$ws.Range("H136").Value = 1985158
$ws.Range("I136").Value = 2268598
$ws.Range("J136").Value = 1077.4286
$ws.Range("K136").Value = 6805794
$ws.Range("L136").Value = 3232.2858
$ws.Range("M136").Value = -6803244
$ws.Range("N136").Value = -8332.2858
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14244.444
$ws.Range("J82").Value = 21559.555
$ws.Range("L82").Value = 21559.555
$ws.Range("N82").Value = -22325.555
$ws.Range("H85").Value = 14244.444
$ws.Range("J85").Value = 21559.555
$ws.Range("L85").Value = 21559.555
$ws.Range("N85").Value = -24211.555
$ws.Range("H134").Value = 3142999
$ws.Range("I134").Value = 4022759.2
$ws.Range("J134").Value = 998.0714
$ws.Range("K134").Value = 12068277.6
$ws.Range("L134").Value = 2994.2142
$ws.Range("M134").Value = -12065742.6
$ws.Range("N134").Value = -8064.2142
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2797687.8
$ws.Range("I31").Value = 997.4693600000001
$ws.Range("J31").Value = 10858736
$ws.Range("K31").Value = 997.4693600000001
$ws.Range("L31").Value = 10858736
$ws.Range("M31").Value = -702.4693600000001
$ws.Range("N31").Value = -10859326
$ws.Range("H34").Value = 2797687.8
$ws.Range("I34").Value = 997.4693600000001
$ws.Range("J34").Value = 10858736
$ws.Range("K34").Value = 997.4693600000001
$ws.Range("L34").Value = 10858736
$ws.Range("M34").Value = -795.4693600000001
$ws.Range("N34").Value = -10859140
$ws.Range("H50").Value = 10127.223
$ws.Range("J50").Value = 11300.333
$ws.Range("L50").Value = 11300.333
$ws.Range("N50").Value = -12550.333
$ws.Range("H51").Value = 13942.556
$ws.Range("J51").Value = 14357.571
$ws.Range("L51").Value = 14357.571
$ws.Range("N51").Value = -15829.571
$ws.Range("H58").Value = 1785.1794
$ws.Range("I58").Value = 962.5319
$ws.Range("J58").Value = 3032.4194
$ws.Range("K58").Value = 962.5319
$ws.Range("L58").Value = 3032.4194
$ws.Range("M58").Value = -759.5319
$ws.Range("N58").Value = -3438.4194
$ws.Range("H60").Value = 8575.77
$ws.Range("J60").Value = 10056.223
$ws.Range("L60").Value = 10056.223
$ws.Range("N60").Value = -11078.223
$ws.Range("H61").Value = 13942.556
$ws.Range("J61").Value = 14357.571
$ws.Range("L61").Value = 14357.571
$ws.Range("N61").Value = -15053.571
$ws.Range("H68").Value = 16541.916
$ws.Range("J68").Value = 18722.555
$ws.Range("L68").Value = 18722.555
$ws.Range("N68").Value = -20220.555
$ws.Range("H71").Value = 16541.916
$ws.Range("J71").Value = 18722.555
$ws.Range("L71").Value = 56167.665
$ws.Range("N71").Value = -63655.665
$ws.Range("H132").Value = 1348.2858
$ws.Range("I132").Value = 1308.3334
$ws.Range("K132").Value = 3925.0002
$ws.Range("M132").Value = -1395.0002
$ws.Range("H134").Value = 1029.7843
$ws.Range("I134").Value = 1064.1628
$ws.Range("J134").Value = 845
$ws.Range("K134").Value = 3192.4884
$ws.Range("L134").Value = 2535
$ws.Range("M134").Value = -657.4884000000002
$ws.Range("N134").Value = -7605
$ws.Range("H136").Value = 1785.1794
$ws.Range("I136").Value = 962.5319
$ws.Range("J136").Value = 3032.4194
$ws.Range("K136").Value = 2887.5957
$ws.Range("L136").Value = 9097.2582
$ws.Range("M136").Value = -337.5956999999999
$ws.Range("N136").Value = -14197.2582
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7813146.5
$ws.Range("I5").Value = 292.96667
$ws.Range("J5").Value = 14706841
$ws.Range("K5").Value = 878.9000100000001
$ws.Range("L5").Value = 44120523
$ws.Range("M5").Value = -766.9000100000001
$ws.Range("N5").Value = -44120747
$ws.Range("H135").Value = 7813146.5
$ws.Range("I135").Value = 292.96667
$ws.Range("J135").Value = 14706841
$ws.Range("K135").Value = 2636.70003
$ws.Range("L135").Value = 132361569
$ws.Range("M135").Value = -101.70003
$ws.Range("N135").Value = -132366639
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27028964
$ws.Range("I132").Value = 37039004
$ws.Range("J132").Value = 1862.8
$ws.Range("K132").Value = 111117012
$ws.Range("L132").Value = 5588.4
$ws.Range("M132").Value = -111114482
$ws.Range("N132").Value = -10648.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1772.8
$ws.Range("I61").Value = 1670.6666
$ws.Range("K61").Value = 1670.6666
$ws.Range("M61").Value = -1468.6666
$ws.Range("H113").Value = 1772.8
$ws.Range("I113").Value = 1670.6666
$ws.Range("K113").Value = 1670.6666
$ws.Range("M113").Value = 499.3334
$ws.Range("H132").Value = 3266.7083
$ws.Range("I132").Value = 3266.7083
$ws.Range("K132").Value = 9800.124899999999
$ws.Range("M132").Value = -7270.124899999999
$ws.Range("H136").Value = 962.5897
$ws.Range("I136").Value = 550.1818
$ws.Range("J136").Value = 3230.8333
$ws.Range("K136").Value = 1650.5454
$ws.Range("L136").Value = 9692.499899999999
$ws.Range("M136").Value = 899.4546
$ws.Range("N136").Value = -14792.4999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 31950
$ws.Range("J114").Value = 31950
$ws.Range("L114").Value = 31950
$ws.Range("N114").Value = -40628
$ws.Range("H132").Value = 6227629.5
$ws.Range("I132").Value = 6538995.5
$ws.Range("J132").Value = 316.66666
$ws.Range("K132").Value = 19616986.5
$ws.Range("L132").Value = 949.9999799999999
$ws.Range("M132").Value = -19614456.5
$ws.Range("N132").Value = -6009.99998
